# Updates crypto price/volume figures per the Jan 1 2023 symbol-list refresh.
# Cells are plain text (inline strings) in the source sheet, so each write
# forces Text format first and restores the default "Normal" style
# afterwards to avoid leaving a stray NumberFormat behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "244.02"
Set-TextValue "E2" "-0.51%"
Set-TextValue "D3" "26.45"
Set-TextValue "E3" "3.48%"
Set-TextValue "D4" "5.135"
Set-TextValue "E4" "0.10%"
Set-TextValue "D5" "0.05606"
Set-TextValue "D6" "6.469"
Set-TextValue "E6" "-0.17%"
Set-TextValue "D7" "0.8188"
Set-TextValue "E7" "0.13%"
Set-TextValue "D8" "0.8333"
Set-TextValue "E8" "-0.93%"
Set-TextValue "D9" "0.1332"
Set-TextValue "E9" "-0.09%"
Set-TextValue "D10" "0.06936"
Set-TextValue "E10" "-0.33%"
Set-TextValue "E11" "1.12%"
Set-TextValue "D12" "0.09379"
Set-TextValue "E12" "-0.04%"
Set-TextValue "D13" "0.001523"
Set-TextValue "E13" "0.17%"
Set-TextValue "D14" "0.0005961"
Set-TextValue "E14" "-0.65%"
Set-TextValue "D15" "0.006155"
Set-TextValue "E15" "-0.92%"
Set-TextValue "D16" "3.653"
Set-TextValue "E16" "3.54%"
Set-TextValue "D17" "3.024"
Set-TextValue "E17" "0.23%"
Set-TextValue "D18" "2.301"
Set-TextValue "D20" "0.03090"
Set-TextValue "E20" "-3.72%"
Set-TextValue "E21" "-1.49%"
Set-TextValue "D22" "3.744"
Set-TextValue "E22" "-0.23%"
Set-TextValue "D23" "0.04592"
Set-TextValue "E23" "-2.31%"
Set-TextValue "E24" "-2.33%"
Set-TextValue "D25" "0.001226"
Set-TextValue "E25" "-1.64%"
Set-TextValue "D26" "0.004494"
Set-TextValue "E26" "-2.48%"
Set-TextValue "D27" "0.00009601"
Set-TextValue "E27" "-1.02%"
Set-TextValue "D28" "0.0001399"
Set-TextValue "E28" "0.69%"
Set-TextValue "D40" "0.03641"
Set-TextValue "E40" "-0.33%"
Set-TextValue "D41" "0.006171"
Set-TextValue "E41" "82.96%"
Set-TextValue "E42" "-22.41%"
Set-TextValue "D43" "0.002400"
Set-TextValue "E43" "-4.66%"
Set-TextValue "D44" "0.008113"
Set-TextValue "E44" "6.13%"
Set-TextValue "D45" "0.00005353"
Set-TextValue "E46" "0.01%"
Set-TextValue "D47" "0.1400"
Set-TextValue "E47" "4.88%"
Set-TextValue "D48" "0.002465"
Set-TextValue "E48" "16.06%"
Set-TextValue "D49" "0.00002100"
Set-TextValue "E49" "0.01%"
Set-TextValue "D50" "0.0002000"
Set-TextValue "E50" "0.01%"
